# Add a "Deny" button/column (G) to the next-possible-queues rules table,
# mirroring the existing "Default return queue" (F) column, and fix the
# F17 header label which was mistakenly left reading "Default next queue".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: extend the "true" sequential flag into G11 (copy keeps the
#     string type + fill formatting instead of turning into a boolean) ---
$ws.Range("F11").Copy($ws.Range("G11"))

# --- Row 12: make sure the new column doesn't carry over the old
#     leftover style that is being retired ---
$ws.Range("G12").Style = "Normal"

# --- Row 13: RuleTable header row - blank cell under new column ---
$ws.Range("F13").Copy($ws.Range("G13"))

# --- Row 14: ACTION header repeated across columns ---
$ws.Range("F14").Copy($ws.Range("G14"))

# --- Row 15: blank spacer row ---
$ws.Range("F15").Copy($ws.Range("G15"))

# --- Row 16: action snippet row - add the new "deny queue" action ---
$ws.Range("F16").Copy($ws.Range("G16"))
$ws.Range("G16").Value = '$model.setDefaultDenyQueue($param);'

# --- Row 17: column headers - fix F17 label and add G17 label ---
$ws.Range("F17").Copy($ws.Range("G17"))
$ws.Range("G17").Value = "Default deny queue"
$ws.Range("F17").Value = "Default return queue"

# --- Row 18: default/null values row ---
$ws.Range("F18").Copy($ws.Range("G18"))

# --- Row heights: widen rows that now wrap more content ---
$ws.Rows.Item(10).RowHeight = 195
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75
$ws.Rows.Item(17).RowHeight = 150
$ws.Rows.Item(18).RowHeight = 30

# --- Leave the selection where the author ended up after adding the
#     new column ---
$ws.Activate()
$ws.Range("F17").Select() | Out-Null
